$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AUC / CI Lower Bound / CI Upper Bound values (new AzBio risk score results)
$ws.Range("C2").Value = 0.6858660130718954
$ws.Range("D2").Value = 0.6211878253438364
$ws.Range("E2").Value = 0.7468058149296148

$ws.Range("C3").Value = 0.6376633986928105
$ws.Range("D3").Value = 0.5794586952998457
$ws.Range("E3").Value = 0.6961786627229969

$ws.Range("C4").Value = 0.6993464052287582
$ws.Range("D4").Value = 0.6333323994747817
$ws.Range("E4").Value = 0.7637974539586546
